$wb = $excel.ActiveWorkbook
$wsQ = $wb.Worksheets.Item("QuickStart")

# --- Set cell values for new rows 28-38, in the same order the original
#     author entered them (so new shared-string indices line up) ---
$wsQ.Range("A28").Value = 'First Groovy program'
$wsQ.Range("B28").Value = '1. Open any folder like E:\practiceProjects\Groovy
2. create Sample.groovy file
3. write println "Welcome to Groovy" in that file
4. open command prompt and go to location where Sample.groovy file is located
5. run groovyc Sample.groovy -> Now Sample.class file will be generated
6. run groovy Sample -> we can see Welcome to Groovy on command prompt'
$wsQ.Range("A30").Value = 'print today''s date'
$wsQ.Range("B30").Value = 'println "Today''s date " + new Date()'
$wsQ.Range("B31").Value = 'We don’t need to add import because automatically adds default/commonly used packages to classpath.'
$wsQ.Range("B29").Value = 'We can directly run groovy Sample.groovy -> this will execute groovy and gives output'
$wsQ.Range("A32").Value = 'Groovy in eclipse'
$wsQ.Range("B32").Value = '1. Install Groovy plugin
2. Eclipse Help -> Market place -> '
$wsQ.Range("A33").Value = 'Create a variable sum'
$wsQ.Range("B33").Value = 'def sum = 3 + 5'
$wsQ.Range("B34").Value = 'def name = "Johny"
println "$name"'
$wsQ.Range("A34").Value = 'declare a string and use it any where - this $ syntax is called G String'
$wsQ.Range("A35").Value = 'Collection iteration with G String'
$wsQ.Range("B35").Value = 'def beatles = ["johny", "depp", "jack"]
def hello = "Hello, "
for(int i=0;i<beatles.size();i++){
 println "$hello" + beatles[i]
}'
$wsQ.Range("B36").Value = 'for(int i=0;i<beatles.size();i++){
 println "$hello" + "${beatles[i]}"
}'
$wsQ.Range("B37").Value = 'for(def i=0;i<beatles.size;i++){
 println "$hello" + "${beatles[i]}"
}'
$wsQ.Range("A38").Value = 'iterate collection with for each loop'
$wsQ.Range("B38").Value = 'for(beatle in beatles){
 println "$hello" + "$beatle"
}'

# --- Apply "wrap text" style (same as existing B14, style index 5) to column B cells ---
#     (PasteSpecial is applied one cell at a time: this engine only honours the
#      first area of a multi-area destination range.)
$wsQ.Range("B14").Copy()
foreach ($r in 28..38) {
    $wsQ.Range("B$r").PasteSpecial(-4122)
}

# --- Create / apply the distinct "reset" style (index 9) used for the A-column
#     cells that participate in vertical merges in the new block ---
$wsQ.Range("A28").WrapText = $False
$wsQ.Range("A28").Copy()
foreach ($addr in @("A29","A30","A31","A35","A36","A37")) {
    $wsQ.Range($addr).PasteSpecial(-4122)
}

# --- Row heights ---
$wsQ.Rows.Item(28).RowHeight = 90
$wsQ.Rows.Item(32).RowHeight = 30
$wsQ.Rows.Item(34).RowHeight = 30
$wsQ.Rows.Item(35).RowHeight = 90
$wsQ.Rows.Item(36).RowHeight = 45
$wsQ.Rows.Item(37).RowHeight = 45
$wsQ.Rows.Item(38).RowHeight = 45

# --- Merge cells ---
$wsQ.Range("A28:A29").Merge()
$wsQ.Range("A30:A31").Merge()
$wsQ.Range("A35:A37").Merge()

# --- Update selection / active cell to mirror final view state ---
$wsQ.Range("A38").Select()

